# Updated symbol list on Sat Dec 24 18:53:01 UTC 2022 with GitHub Actions
#
# Re-applies the refreshed coin-ranking snapshot: a handful of standalone
# price refreshes, a block of rows (9-17) that rotate up by one position
# (the "One" coin drops out of the #8 slot and reappears at #16, bumping
# every row between up by one), and a couple of trailing "Bestin24h" /
# "Worstin24h" label tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-like cells are stored as plain text in this workbook (t="inlineStr").
# Assigning a numeric-looking string straight to .Value makes Excel coerce
# it into a real number, which would change both the stored type and the
# precision (e.g. "5.385" -> 5.3849999999999998). Force each cell to a text
# number-format first so the literal string is preserved, then drop the
# format override again so no stray style sticks around on the cell.
function Set-TextValue {
    param($Address, $Value)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

# --- Standalone price refreshes -------------------------------------------
Set-TextValue "D3" "21.92"
Set-TextValue "D4" "5.385"
Set-TextValue "D7" "0.8154"
Set-TextValue "D8" "0.9389"

# --- Rows 9-17: "One" rotates from rank 8 to rank 16, shifting the coins
#     that used to occupy ranks 9-16 each up one row, with a few prices
#     refreshed along the way. --------------------------------------------
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1427"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07435"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03523"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03057"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09417"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.006"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001600"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04824"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005942"
$ws.Range("E17").Value = "16OneONE"

# --- Remaining standalone price refreshes ----------------------------------
Set-TextValue "D18" "0.005203"
Set-TextValue "D19" "0.004159"
Set-TextValue "D20" "0.0009920"
Set-TextValue "D21" "3.667"
Set-TextValue "D22" "6.417"
Set-TextValue "D26" "0.00007000"
Set-TextValue "D40" "0.04006"
Set-TextValue "D41" "0.006450"
Set-TextValue "D42" "0.1075"
Set-TextValue "D43" "0.002720"
Set-TextValue "D44" "0.005931"
Set-TextValue "D45" "0.00005294"
Set-TextValue "D48" "0.002395"
Set-TextValue "D49" "0.00002101"

# --- "Best/Worst in 24h" badge changes --------------------------------------
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
